# Scenario7-Chaithanyaprod-TR2-LATE-Makerepayment2.xlsx
# "chaitanya basic scenarios add"
#
# This script:
#   1) Updates literal data values on the "Summary" and "Transactions" sheets.
#   2) Re-points each sheet's saved cursor/selection (and, where applicable,
#      the scrolled top-left cell) to match the values captured when the
#      workbook was last saved by the author.
#   3) Leaves "Input" as the active (tabSelected) sheet/selection, since that
#      is the sheet that was showing when the workbook was saved.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Data edits (no UI/selection side effects)
# ---------------------------------------------------------------------------

# Summary!F3 : 4171.83 -> 4171.82
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("F3").Value = 4171.82

# Transactions!A2:A5 : Entry ID column renumbered
$wsTransactions = $wb.Worksheets.Item("Transactions")
$wsTransactions.Range("A2").Value = 236
$wsTransactions.Range("A3").Value = 231
$wsTransactions.Range("A4").Value = 229
$wsTransactions.Range("A5").Value = 227

# ---------------------------------------------------------------------------
# View-state edits (selection / scroll position per sheet)
# ---------------------------------------------------------------------------

# Summary : selection C21 -> C4
$wsSummary.Range("C4").Select()

# Transactions : selection D4 -> C4
$wsTransactions.Range("C4").Select()

# Repayment Schedule : scrolled so row 7 is the top visible row, selection C12 -> C11
$wsRepayment = $wb.Worksheets.Item("Repayment Schedule")
$wsRepayment.Activate()
$excel.ActiveWindow.ScrollRow = 7
$wsRepayment.Range("C11").Select()

# Input : selection B3 -> A2 (kept as the active sheet/tab, selected last)
$wsInput = $wb.Worksheets.Item("Input")
$wsInput.Range("A2").Select()
